$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two extra data rows (old rows 3 and 4), leaving header + 1 data row
$ws.Rows.Item(4).Delete() | Out-Null
$ws.Rows.Item(3).Delete() | Out-Null

# Add a 4th header column, copying the existing header style (bold/border/center)
$ws.Range("C1").Copy($ws.Range("D1")) | Out-Null
$ws.Range("D1").Value = "Atualizacao"

# Re-label the existing headers (columns shifted meaning)
$ws.Range("B1").Value = "Mapa Selecionado"
$ws.Range("C1").Value = "Numero endereco"

# Update the remaining data row
$ws.Range("A2").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("A2").Value = [DateTime]"2024-10-05"

$ws.Range("B2").Value = "mapa1"

# "3" must be stored as text, not a number, matching the source data
$ws.Range("C2").Value = "'3"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").Value = "mudou-se"

# Column widths to match the new layout (nearest value Excel's pixel/MDW
# quantization will round to the target stored width)
$ws.Columns.Item(1).ColumnWidth = 10.62
$ws.Columns.Item(2).ColumnWidth = 15.77
$ws.Columns.Item(3).ColumnWidth = 15.44
$ws.Columns.Item(4).ColumnWidth = 9.77

$ws.Range("C12").Select() | Out-Null
